$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.877531886100769
$ws.Range("B1").Value = 2.073564529418945
$ws.Range("C1").Value = 2.456933736801147
$ws.Range("D1").Value = 3.785820484161377
$ws.Range("E1").Value = 1.14708685874939
